$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old sub-header row (row 2: "(MW)/(MW)/(GWh)/(GWh)/(GWh)") - shifts data rows up
$ws.Rows(2).Delete()

# Rewrite the header row with the new column titles
$ws.Range("A1").Value2 = "idx"
$ws.Range("B1").Value2 = "idx2"
$ws.Range("C1").Value2 = "Name"
$ws.Range("D1").Value2 = "Date Start"
$ws.Range("E1").Value2 = "Date End"
$ws.Range("F1").Value2 = "(m3/s)"
$ws.Range("G1").Value2 = "(MW1)"
$ws.Range("H1").Value2 = "(MW2)"
$ws.Range("I1").Value2 = "(GWh) Winter"
$ws.Range("J1").Value2 = "(GWh) Summer"
$ws.Range("K1").Value2 = "(GWh) Year"

# Update selection to the first data row
$ws.Range("A2:K2").Select()
